# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets
# to reflect newly generated output (gh-pages rebuild at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1394
$ws1.Range("F5").Value = 112
$ws1.Range("F6").Value = 67
$ws1.Range("F7").Value = 11832
$ws1.Range("F8").Value = 4423
$ws1.Range("F9").Value = 32
$ws1.Range("F10").Value = 46
$ws1.Range("F14").Value = 1104
$ws1.Range("F15").Value = 159
$ws1.Range("F17").Value = 5135
$ws1.Range("F21").Value = 11375
$ws1.Range("F22").Value = 11341

# --- Sheet "全部类型" (all types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1394
$ws4.Range("F5").Value = 112
$ws4.Range("F6").Value = 67
$ws4.Range("F7").Value = 11832
$ws4.Range("F8").Value = 4423
$ws4.Range("F9").Value = 32
$ws4.Range("F10").Value = 46
$ws4.Range("F15").Value = 1104
$ws4.Range("F16").Value = 159
$ws4.Range("F18").Value = 5135
$ws4.Range("F22").Value = 11375
$ws4.Range("F23").Value = 11341

$wb.Save()
